$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.616.31"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "1.878.60"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.52"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4765"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2930"
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06533"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7412"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.76"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "1.876.18"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.215"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.00"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "30.715.65"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.26"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007534"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "2.123.32"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.265"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.208"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.45"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.201"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.90"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.917"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09851"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.340"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.293"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.119"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04834"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6962"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.762"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.271"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.44"
$ws.Range("E41").Value = "  +5.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.994"
$ws.Range("E42").Value = "  +4.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4245"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8384"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.27"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.368"
$ws.Range("E47").Value = "  +1.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.033"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "910.57"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3919"
$ws.Range("E51").Value = "  +2.26%  "
